# Applies the cryptos.xlsx price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of "Price" column values look like plain decimal numbers (single
# decimal point, e.g. "243.11"). The source workbook stores every Price/Volume
# cell as text, so force a Text number format on just those cells first to stop
# Excel from re-interpreting the literal string as a numeric value. (Looping
# cell-by-cell because applying NumberFormat to a multi-area union Range only
# affects the first area.)
$textForceCells = @("D5","D8","D10","D11","D12","D13","D15","D16","D18","D20","D22","D23","D26","D27","D29","D30","D35","D38","D39","D41","D43","D44","D45","D46","D47","D48","D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '37.015.00'
$ws.Range('E2').Value = '  -1.01%  '

$ws.Range('D3').Value = '1.997.00'
$ws.Range('E3').Value = '  -1.35%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '243.11'
$ws.Range('E5').Value = '  -4.18%  '

$ws.Range('E6').Value = '  -2.32%  '

$ws.Range('E7').Value = '  +0.09%  '

$ws.Range('D8').Value = '54.42'
$ws.Range('E8').Value = '  -4.44%  '

$ws.Range('E9').Value = '  -2.44%  '

$ws.Range('D10').Value = '57.14'
$ws.Range('E10').Value = '  +0.00%  '

$ws.Range('D11').Value = '0.0754'
$ws.Range('E11').Value = '  -4.24%  '

$ws.Range('D12').Value = '0.0981'
$ws.Range('E12').Value = '  -3.43%  '

$ws.Range('D13').Value = '14.20'
$ws.Range('E13').Value = '  -3.53%  '

$ws.Range('D14').Value = '2.289.36'
$ws.Range('E14').Value = '  -1.45%  '

$ws.Range('D15').Value = '21.01'
$ws.Range('E15').Value = '  -0.61%  '

$ws.Range('D16').Value = '0.759'
$ws.Range('E16').Value = '  -6.91%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.019.13'
$ws.Range('E17').Value = '  -0.79%  '

$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '5.07'
$ws.Range('E18').Value = '  -4.85%  '

$ws.Range('D19').Value = '36.917.43'
$ws.Range('E19').Value = '  -0.87%  '

$ws.Range('D20').Value = '68.62'
$ws.Range('E20').Value = '  -1.08%  '

$ws.Range('E21').Value = '  -4.29%  '

$ws.Range('D22').Value = '5.12'

$ws.Range('D23').Value = '228.58'
$ws.Range('E23').Value = '  +0.18%  '

$ws.Range('E24').Value = '  +0.20%  '

$ws.Range('E25').Value = '  -6.49%  '

$ws.Range('D26').Value = '2.35'
$ws.Range('E26').Value = '  +0.54%  '

$ws.Range('D27').Value = '162.46'
$ws.Range('E27').Value = '  -0.07%  '

$ws.Range('E28').Value = '  -3.98%  '

$ws.Range('D29').Value = '19.22'
$ws.Range('E29').Value = '  -3.15%  '

$ws.Range('D30').Value = '0.126'
$ws.Range('E30').Value = '  -4.09%  '

$ws.Range('E31').Value = '  -3.34%  '

$ws.Range('E32').Value = '  -1.73%  '

$ws.Range('E33').Value = '  -4.99%  '

$ws.Range('E34').Value = '  -7.00%  '

$ws.Range('D35').Value = '4.25'
$ws.Range('E35').Value = '  -5.71%  '

$ws.Range('E36').Value = '  -5.89%  '

$ws.Range('E37').Value = '  -0.10%  '

$ws.Range('D38').Value = '1.79'
$ws.Range('E38').Value = '  -1.88%  '

$ws.Range('D39').Value = '3.26'
$ws.Range('E39').Value = '  -4.77%  '

$ws.Range('E40').Value = '  +0.57%  '

$ws.Range('D41').Value = '3.07'
$ws.Range('E41').Value = '  +1.35%  '

$ws.Range('D42').Value = '1.429.44'
$ws.Range('E42').Value = '  +2.01%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0204'
$ws.Range('E43').Value = '  -4.80%  '

$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = '0.0895'
$ws.Range('E44').Value = '  -7.00%  '

$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '1.13'
$ws.Range('E45').Value = '  -4.39%  '

$ws.Range('D46').Value = '88.35'
$ws.Range('E46').Value = '  -1.93%  '

$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '1.01'
$ws.Range('E47').Value = '  -2.40%  '

$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '15.17'
$ws.Range('E48').Value = '  -5.34%  '

$ws.Range('E49').Value = '  +0.62%  '

$ws.Range('D50').Value = '6.73'
$ws.Range('E50').Value = '  -8.15%  '

$ws.Range('D51').Value = '2.180.85'
$ws.Range('E51').Value = '  -1.45%  '
